$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Translations - Lab")

$rows = @(
    @("cs", "lab.inventory.atomizer.index.title", "Atomizéry"),
    @("cs", "lab.inventory.mod.index.title", "Mody"),
    @("cs", "lab.inventory.cell.index.title", "Články"),
    @("cs", "lab.inventory.cotton.index.title", "Vaty"),
    @("cs", "lab.inventory.aroma.index.title", "Aromata"),
    @("cs", "lab.inventory.base.index.title", "Báze"),
    @("cs", "lab.inventory.booster.index.title", "Boostery")
)

$startRow = 18
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
    $ws.Cells.Item($r, 3).Value = $rows[$i][2]
}

$ws.Range("A18:C24").WrapText = $true
$ws.Range("A18:C24").Font.Size = 10

$ws.Range("B19").Select()
